$d = $word.ActiveDocument
$d.Content.Find.Execute(" : T1.1", $true, $false, $false, $false, $false,
                         $true, 1, $false, " : T1", 2)
